$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 224 (pushing old rows 224..347 down to 226..349)
$ws.Rows.Item(224).Insert()
$ws.Rows.Item(224).Insert()

# --- New row 224 ---
# Same constant metadata columns as every other data row on this sheet,
# but with its own date/price/volume data.
$ws.Cells.Item(224, 1).Value = 4
$ws.Cells.Item(224, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(224, 3).Value = "Los Lagos"
$ws.Cells.Item(224, 4).Value = 44777
$ws.Cells.Item(224, 5).Value = 10
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100102
$ws.Cells.Item(224, 8).Value = "Cítricos"
$ws.Cells.Item(224, 9).Value = 100102006
$ws.Cells.Item(224, 10).Value = "Pomelo"
$ws.Cells.Item(224, 11).Value = "Start Ruby"
$ws.Cells.Item(224, 12).Value = "Primera"
$ws.Cells.Item(224, 13).Value = 120
$ws.Cells.Item(224, 14).Value = 14000
$ws.Cells.Item(224, 15).Value = 15000
$ws.Cells.Item(224, 16).Value = 14500
$ws.Cells.Item(224, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(224, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(224, 19).Value = 1036
$ws.Cells.Item(224, 20).Value = 14

# --- New row 225 ---
$ws.Cells.Item(225, 1).Value = 4
$ws.Cells.Item(225, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(225, 3).Value = "Los Lagos"
$ws.Cells.Item(225, 4).Value = 44777
$ws.Cells.Item(225, 5).Value = 10
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100102
$ws.Cells.Item(225, 8).Value = "Cítricos"
$ws.Cells.Item(225, 9).Value = 100102006
$ws.Cells.Item(225, 10).Value = "Pomelo"
$ws.Cells.Item(225, 11).Value = "Start Ruby"
$ws.Cells.Item(225, 12).Value = "Segunda"
$ws.Cells.Item(225, 13).Value = 60
$ws.Cells.Item(225, 14).Value = 12000
$ws.Cells.Item(225, 15).Value = 12000
$ws.Cells.Item(225, 16).Value = 12000
$ws.Cells.Item(225, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(225, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(225, 19).Value = 857
$ws.Cells.Item(225, 20).Value = 14
